# C3DC Regression/Smoke suite update (TC01_C3DC_phs002677_SexAtBirth-Female.xlsx)
#
# The "TreatmentTab" row's query (cell B5) wrapped the REPLACE(...) call in a
# redundant CONCAT(...) — CONCAT of a single argument is a no-op, so it is
# removed, leaving a plain REPLACE(...) call for the "Treatment Agent" column.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$oldQuery = $treatmentCell.Value2
$newQuery = $oldQuery.Replace( `
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))", `
    "REPLACE(trt.treatment_agent, ';', ', ')")

$treatmentCell.Value2 = $newQuery

# Leave the workbook with cell B2 as the active selection, matching the
# saved view state from the authoring session.
$ws.Range("B2").Select()
